{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// The document currently has a single (empty) paragraph. Insert the new\n// \"Difference between:\" heading paragraph before it, and the four new\n// code/comment paragraphs after it, leaving the original empty paragraph\n// untouched in between.\nconst firstParagraph = paragraphs.items[0];\n\nfirstParagraph.insertParagraph(\"Difference between:\", \"Before\");\n\nlet cursor = firstParagraph;\ncursor = cursor.insertParagraph(\n  \"# Get the length-frequency data (length in millimeters)\",\n  \"After\"\n);\ncursor = cursor.insertParagraph(\n  'lfmm <- read.lfmmdata.f(raw_data_dir,\"LengthMM2000-2021.txt\")',\n  \"After\"\n);\ncursor = cursor.insertParagraph(\n  \"# Get the grouped length-frequency output\",\n  \"After\"\n);\ncursor = cursor.insertParagraph(\n  'lfgrpd <- read.lengthfreq.f(raw_data_dir,\"LengthFreq2000-2021.txt\")',\n  \"After\"\n);\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# The document currently has a single empty paragraph. Insert the new\n# \"Difference between:\" heading paragraph before it.\n$firstPara = $d.Paragraphs(1)\n$firstPara.Range.InsertParagraphBefore()\n$d.Paragraphs(1).Range.Text = \"Difference between:\"\n\n# The original empty paragraph is now Paragraphs(2); insert the four new\n# code/comment paragraphs after it, one at a time, so they land in order.\n$anchor = $d.Paragraphs(2)\n$anchor.Range.InsertParagraphAfter()\n$d.Paragraphs(3).Range.Text = \"# Get the length-frequency data (length in millimeters)\"\n\n$anchor = $d.Paragraphs(3)\n$anchor.Range.InsertParagraphAfter()\n$d.Paragraphs(4).Range.Text = 'lfmm <- read.lfmmdata.f(raw_data_dir,\"LengthMM2000-2021.txt\")'\n\n$anchor = $d.Paragraphs(4)\n$anchor.Range.InsertParagraphAfter()\n$d.Paragraphs(5).Range.Text = \"# Get the grouped length-frequency output\"\n\n$anchor = $d.Paragraphs(5)\n$anchor.Range.InsertParagraphAfter()\n$d.Paragraphs(6).Range.Text = 'lfgrpd <- read.lengthfreq.f(raw_data_dir,\"LengthFreq2000-2021.txt\")'\n"}
